# Update handback timestamp values on the "Generate Report for Handback" run.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 0e8879ea... row (row 3)
$wsOverview.Range("G3").Value = "2016-09-02 18:50:44"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 3
$wsZhCn.Range("H3").Value = "2016-09-02 18:50:39"
$wsZhCn.Range("K3").Value = "2016-09-02 18:50:57"

# de-de sheet: Correspond Handback DateTime for row 3
$wsDeDe.Range("K3").Value = "2016-09-02 18:51:12"
